$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet name to reflect the new "through" date
$ws.Name = "Through 2022-10-29"

# Update the header label in I1 (shared string "2022 (through 10-28)")
$ws.Range("I1").Value = "2022 (through 10-29)"

# Update data values for new day (2022-11-06)
$ws.Range("I2").Value = 162
$ws.Range("I11").Value = 114
$ws.Range("I14").Value = 1391
